# Fruta / hortaliza, semanal
# Insert a new price record as row 30 in the daily logic sheet for
# "Vega Modelo de Temuco - Arándano (blue)". This pushes the existing
# rows 30..88 down by one (to 31..89) and fills the newly inserted row
# with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 30, shifting rows 30-88 down to 31-89.
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with the new record.
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 44662
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100101
$ws.Range("H30").Value = "Berries"
$ws.Range("I30").Value = 100101001
$ws.Range("J30").Value = "Arándano (blue)"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = 2000
$ws.Range("O30").Value = 2000
$ws.Range("P30").Value = 2000
$ws.Range("Q30").Value = "$/kilo"
$ws.Range("R30").Value = "Región de La Araucanía"
$ws.Range("S30").Value = 2000
$ws.Range("T30").Value = 1
